$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 31   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/10/2024  Through  6/16/2024"

# --- Crime statistics table updates (rows 15-28) ---
$ws.Range("C14").Copy($ws.Range("D15"))
$ws.Range("E14").Copy($ws.Range("E15"))
$ws.Range("M15").Value = 11.111111111111
$ws.Range("N15").Value = -33.333333333333
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("C14").Copy($ws.Range("D16"))
$ws.Range("E14").Copy($ws.Range("E16"))
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 20
$ws.Range("L16").Value = 30.555555555555
$ws.Range("M16").Value = -16.071428571428
$ws.Range("N16").Value = -82.783882783882
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 28
$ws.Range("H17").Value = 21.739130434782
$ws.Range("I17").Value = 151
$ws.Range("J17").Value = 106
$ws.Range("K17").Value = 42.452830188679
$ws.Range("L17").Value = 48.039215686274
$ws.Range("M17").Value = 109.722222222222
$ws.Range("N17").Value = -23.350253807106
$ws.Range("C18").Value = 1
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = -11.111111111111
$ws.Range("L18").Value = 3.225806451612
$ws.Range("M18").Value = -41.818181818181
$ws.Range("N18").Value = -90.184049079754
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -66.666666666666
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -17.647058823529
$ws.Range("I19").Value = 72
$ws.Range("J19").Value = 84
$ws.Range("K19").Value = -14.285714285714
$ws.Range("L19").Value = -20
$ws.Range("M19").Value = 56.521739130434
$ws.Range("N19").Value = -42.4
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -90
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -22.857142857142
$ws.Range("N20").Value = -88.559322033898
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -31.25
$ws.Range("F21").Value = 53
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = -15.873015873015
$ws.Range("I21").Value = 341
$ws.Range("J21").Value = 320
$ws.Range("K21").Value = 6.5625
$ws.Range("L21").Value = 16.780821917808
$ws.Range("M21").Value = 22.661870503597
$ws.Range("N21").Value = -71.101694915254
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("M22").Value = 20
$ws.Range("C23").Value = 4
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 11.111111111111
$ws.Range("I23").Value = 43
$ws.Range("J23").Value = 44
$ws.Range("K23").Value = -2.272727272727
$ws.Range("L23").Value = 22.857142857142
$ws.Range("M23").Value = 126.315789473684
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 53
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = 32.5
$ws.Range("I24").Value = 276
$ws.Range("J24").Value = 256
$ws.Range("K24").Value = 7.8125
$ws.Range("L24").Value = 7.8125
$ws.Range("M24").Value = 64.285714285714
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("C14").Copy($ws.Range("D25"))
$ws.Range("E14").Copy($ws.Range("E25"))
$ws.Range("F25").Value = 3
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = -57.142857142857
$ws.Range("L25").Value = -38.235294117647
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 225
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = 11.764705882352
$ws.Range("I26").Value = 180
$ws.Range("J26").Value = 175
$ws.Range("K26").Value = 2.857142857142
$ws.Range("L26").Value = 9.090909090909
$ws.Range("M26").Value = -5.759162303664
$ws.Range("C27").Value = 1
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 18
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = 10.526315789473
$ws.Range("L28").Value = -8.695652173913

Write-Output "edits applied"
